$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
# A1 "Name" stays the same.
$ws.Range("B1").Value = "WorkflowsRoot"
$ws.Range("C1").Value = "FormFile"
$ws.Range("D1").Value = "ModuleSetupWorkflow"
$ws.Range("E1").Value = "ConfigPath"

# --- Row 2 (Dispatcher - Basic) ---
$ws.Range("B2").Value = ".templates\Dispatchers\Basic"
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = ".templates\Dispatchers\Basic\Data\BasicDispatcherConfig.xlsx"

# --- Row 3 (Dispatcher - Application) ---
$ws.Range("B3").Value = ".templates\Dispatchers\Application"
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = ".templates\Dispatchers\Application\Data\ApplicationDispatcherConfig.xlsx"

# --- Row 4 (Performer - Basic) ---
$ws.Range("B4").Value = ".templates\Performers\Basic"
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = ".templates\Performers\Basic\Data\BasicPerformerConfig.xlsx"

# --- Row 5 (Performer - REFramework) ---
$ws.Range("B5").Value = ".templates\Performers\REFramework"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = ".templates\Performers\REFramework\Data\Config.xlsx"

# --- Row 6 (Reporter - Basic) ---
$ws.Range("B6").Value = ".templates\Reporters\Basic"
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = ".templates\Reporters\Basic\Data\BasicReporterConfig.xlsx"

# --- Column widths ---
# Column B now holds the long ".templates\..." path strings (same content
# that used to live in column C), so it ends up the same width column C
# used to have. Columns C and D are given the same width for visual
# consistency even though they no longer hold data.
$ws.Columns.Item(2).ColumnWidth = 34.43
$ws.Columns.Item(3).ColumnWidth = 34.43
$ws.Columns.Item(4).ColumnWidth = 34.43

# --- Selection ---
$ws.Range("D3").Select()
